$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Create / update styles
# ---------------------------------------------------------------

# "Default" paragraph style
$default = $d.Styles.Add("Default", 1)
$default.ParagraphFormat.TabStops.Add(35.45, 0)
$default.ParagraphFormat.LineSpacingRule = 3
$default.ParagraphFormat.LineSpacing = 13.8
$default.Font.Name = "Calibri"
$default.Font.NameFarEast = "DejaVu Sans"

# Replace "Hyperlink" character style with "InternetLink"
$dpf = $d.Styles.Item("DefaultParagraphFont")
$hyperlinkStyle = $d.Styles.Item("Hyperlink")
$hyperlinkStyle.Delete()

$internetLink = $d.Styles.Add("InternetLink", 2)
$internetLink.BaseStyle = $dpf
$internetLink.NameLocal = "Internet Link"
$internetLink.Font.Color = 16711680
$internetLink.Font.Underline = 1

# "Heading" paragraph style
$heading = $d.Styles.Add("Heading", 1)
$heading.BaseStyle = $default
$heading.ParagraphFormat.KeepWithNext = 1
$heading.ParagraphFormat.SpaceBefore = 12
$heading.ParagraphFormat.SpaceAfter = 6
$heading.Font.Name = "Arial"
$heading.Font.Size = 14
$heading.Font.SizeBi = 14
$heading.Font.NameBi = "DejaVu Sans"

# "Textbody" paragraph style
$textbody = $d.Styles.Add("Textbody", 1)
$textbody.BaseStyle = $default
$textbody.NameLocal = "Text body"
$textbody.ParagraphFormat.SpaceAfter = 6

$heading.NextParagraphStyle = $textbody

# "List" paragraph style
$list = $d.Styles.Add("List", 1)
$list.BaseStyle = $textbody

# "Caption" paragraph style
$caption = $d.Styles.Add("Caption", 1)
$caption.BaseStyle = $default
$caption.NameLocal = "caption"
$caption.ParagraphFormat.NoLineNumber = 1
$caption.ParagraphFormat.SpaceBefore = 6
$caption.ParagraphFormat.SpaceAfter = 6
$caption.Font.Italic = 1
$caption.Font.ItalicBi = 1
$caption.Font.Size = 12
$caption.Font.SizeBi = 12

# "Index" paragraph style
$index = $d.Styles.Add("Index", 1)
$index.BaseStyle = $default
$index.ParagraphFormat.NoLineNumber = 1

# ---------------------------------------------------------------
# 2. Apply the "Default" style to all existing paragraphs
# ---------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Style = $default
}

# Paragraph 2 ("This bit is in bold and italic") lost direct formatting
# when the style was applied - re-apply bold/italic to the text run only
# (not the paragraph mark) to avoid re-creating pPr/rPr.
$p2 = $d.Paragraphs.Item(2)
$p2Text = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$p2Text.Font.Bold = 1
$p2Text.Font.Italic = 1

# ---------------------------------------------------------------
# 3. Paragraph 1: append a "." run
# ---------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1End = $d.Range($p1.Range.End - 1, $p1.Range.End - 1)
$p1End.InsertAfter(".")

# ---------------------------------------------------------------
# 4. Insert the new paragraph after "Back to normal" (paragraph 3)
# ---------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$insertPoint = $d.Range($p4.Range.Start, $p4.Range.Start)
$insertPoint.InsertParagraphBefore()

$p4 = $d.Paragraphs.Item(4)
$p4.Style = $default
$cur = $p4.Range.Start

$r = $d.Range($cur, $cur)
$r.InsertAfter("This contains ")
$cur = $r.End

$r = $d.Range($cur, $cur)
$r.InsertAfter("BOLD")
$r.Font.Bold = 1
$r.Font.BoldBi = 1
$cur = $r.End

$r = $d.Range($cur, $cur)
$r.InsertAfter(", ")
$cur = $r.End

$r = $d.Range($cur, $cur)
$r.InsertAfter("ITALIC")
$r.Font.Italic = 1
$r.Font.ItalicBi = 1
$cur = $r.End

$r = $d.Range($cur, $cur)
$r.InsertAfter(" and ")
$cur = $r.End

$r = $d.Range($cur, $cur)
$r.InsertAfter("BOTH")
$r.Font.Bold = 1
$r.Font.BoldBi = 1
$r.Font.Italic = 1
$r.Font.ItalicBi = 1
$cur = $r.End

$r = $d.Range($cur, $cur)
$r.InsertAfter(", as well as ")
$cur = $r.End

$r = $d.Range($cur, $cur)
$r.InsertAfter("RED")
$r.Font.Color = 128
$cur = $r.End

$r = $d.Range($cur, $cur)
$r.InsertAfter(" and ")
$cur = $r.End

$r = $d.Range($cur, $cur)
$r.InsertAfter("YELLOW")
$r.Font.Color = 65510
$cur = $r.End

$r = $d.Range($cur, $cur)
$r.InsertAfter(" text.")
$cur = $r.End

# ---------------------------------------------------------------
# 5. Hyperlink paragraph (now paragraph 5): remove the second
#    hyperlink, rebuild surrounding text, and re-style the first
#    hyperlink.
# ---------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$p5Rng = $d.Range($p5.Range.Start, $p5.Range.End - 1)
$p5Rng.Text = ""

$cur = $p5Rng.Start

$r = $d.Range($cur, $cur)
$r.InsertAfter("We have a ")
$cur = $r.End

$r = $d.Range($cur, $cur)
$r.InsertAfter("hyperlink")
$hlink = $d.Hyperlinks.Add($r, "http://poi.apache.org/")
$cur = $r.End

$r = $d.Range($cur, $cur)
$r.InsertAfter(" here, and another.")
$cur = $r.End

$firstLink = $d.Hyperlinks.Item(1)
$firstLink.Range.Style = "Internet Link"

# ---------------------------------------------------------------
# 6. Section properties
# ---------------------------------------------------------------
$ps = $d.PageSetup
$ps.HeaderDistance = 36
$ps.FooterDistance = 36
$ps.TextColumns.Spacing = 36

$sec = $d.Sections.Item(1)
$sec.ProtectedForForms = $false

Write-Host "Edit complete"
